$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "italic"
$ws.Range("A2").Font.Italic = $true

$ws.Range("A3").Value = "underline"
$ws.Range("A3").Font.Underline = $true

$ws.Range("A4").Value = "strikethrough"
$ws.Range("A4").Font.Strikethrough = $true
